# Update crypto price/volume table cells per the source diff.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.399.43'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '1.844.78'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.017'
$ws.Range("E4").Value = '  +1.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.65'
$ws.Range("E5").Value = '  +2.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.015'
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4743'
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3701'
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07453'
$ws.Range("E9").Value = '  +1.06%  '
$ws.Range("E10").Value = '  +1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.53'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '1.836.64'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07392'
$ws.Range("E13").Value = '  +4.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.500'
$ws.Range("E14").Value = '  +2.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.41'
$ws.Range("E15").Value = '  +1.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.599'
$ws.Range("E16").Value = '  +1.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.016'
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008854'
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.014'
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.86'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("D21").Value = '27.405.26'
$ws.Range("E21").Value = '  +1.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.347'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").Value = '2.066.40'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.912'
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.56'
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  +1.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.184'
$ws.Range("E28").Value = '  +0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.290'
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.06'
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08973'
$ws.Range("E31").Value = '  +0.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7618'
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.178'
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.572'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.957'
$ws.Range("E35").Value = '  +1.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.016'
$ws.Range("E36").Value = '  +1.56%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.110'
$ws.Range("E37").Value = '  +2.01%  '
$ws.Range("E38").Value = '  +1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01972'
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.013'
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.345'
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.412'
$ws.Range("E42").Value = '  +2.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5362'
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1670'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.553'
$ws.Range("E45").Value = '  +1.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4962'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.49'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.017'
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.93'
$ws.Range("E49").Value = '  +1.05%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.684'
$ws.Range("E50").Value = '  +0.85%  '
$ws.Range("E51").Value = '  +0.87%  '
